# For each "sound credit" paragraph (the ones that hold the freesound.org
# hyperlink for a given sound effect), add a new paragraph right after it
# crediting the author, e.g. "by timgormly".
#
# Each sound-credit paragraph is located by a distinctive substring of its
# text (the freesound.org username embedded in the hyperlink URL) instead
# of a hard-coded paragraph index, and the list is processed back-to-front
# so inserting a paragraph never invalidates the index of a paragraph that
# still needs to be processed.

$d = $word.ActiveDocument

$credits = @(
    @{ Marker = "LittleRobotSoundFactory/sounds/270334"; Author = "LittleRobotSoundFactory" },
    @{ Marker = "LittleRobotSoundFactory/sounds/270528"; Author = "LittleRobotSoundFactory" },
    @{ Marker = "Mrthenoronha";                          Author = "Mrthenoronha" },
    @{ Marker = "sharesynth";                            Author = "sharesynth" },
    @{ Marker = "timgormly";                             Author = "timgormly" }
)

foreach ($credit in $credits) {
    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
        $para = $d.Paragraphs.Item($i)
        if ($para.Range.Text -like "*$($credit.Marker)*") {
            $para.Range.InsertParagraphAfter()
            $d.Paragraphs.Item($i + 1).Range.Text = "by " + $credit.Author
            break
        }
    }
}
